# Add the new time-tracking entries for chapters 10 & 11 (rows 9-12 of the
# "Progress" log table). Dates/times are written as the underlying serial
# numbers so they land exactly on the same values Excel itself would store.
# The "Time" column (D) is a calculated table column, and the "Total time"
# cell (G2) is =SUM(D:D), so both recalculate automatically once the new
# Time start / Time end values are entered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: 2025-11-11, 10:00 - 12:30  (no description yet)
$ws.Range("A9").Value = 45972
$ws.Range("B9").Value = 0.41666666666666669
$ws.Range("C9").Value = 0.52083333333333337

# Row 10: 2025-11-11, 13:30 - 16:00 -> "Finished chapter 10"
$ws.Range("A10").Value = 45972
$ws.Range("B10").Value = 0.5625
$ws.Range("C10").Value = 0.66666666666666663
$ws.Range("E10").Value = "Finished chapter 10"

# Row 11: 2025-11-12, 13:15 - 16:15 -> "Almost finished chapter 11"
$ws.Range("A11").Value = 45973
$ws.Range("B11").Value = 0.55208333333333337
$ws.Range("C11").Value = 0.67708333333333337
$ws.Range("E11").Value = "Almost finished chapter 11"

# Row 12: 2025-11-14, 10:30 - 12:00 -> "Finished chapter 11"
$ws.Range("A12").Value = 45975
$ws.Range("B12").Value = 0.4375
$ws.Range("C12").Value = 0.5
$ws.Range("E12").Value = "Finished chapter 11"

# Leave the view scrolled/selected where the author last left it.
$ws.Range("A5").Select()
$ws.Range("E13").Select()
